# Scheduled runner update: refresh market-price-derived columns (H:N)
# across the eight job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 281.6
$ws.Range("I4").Value = 314.5
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 314.5
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -200.5
$ws.Range("N4").Value = -378

# Row 17
$ws.Range("H17").Value = 706220.6
$ws.Range("J17").Value = 812782.2
$ws.Range("L17").Value = 2438346.6
$ws.Range("N17").Value = -2438682.6

# Row 132
$ws.Range("H132").Value = 161903.95
$ws.Range("I132").Value = 169467.97
$ws.Range("K132").Value = 508403.91
$ws.Range("M132").Value = -505873.91

# Row 135
$ws.Range("H135").Value = 1455.7288
$ws.Range("I135").Value = 1314.174
$ws.Range("J135").Value = 1956.6154
$ws.Range("K135").Value = 11827.566
$ws.Range("L135").Value = 17609.5386
$ws.Range("M135").Value = -9292.565999999999
$ws.Range("N135").Value = -22679.5386

# Row 138
$ws.Range("H138").Value = 1237.31
$ws.Range("I138").Value = 756.44446
$ws.Range("J138").Value = 2056.081
$ws.Range("K138").Value = 2269.33338
$ws.Range("L138").Value = 6168.243
$ws.Range("M138").Value = 2870.66662
$ws.Range("N138").Value = -16448.243

# Row 141
$ws.Range("H141").Value = 2300.611
$ws.Range("I141").Value = 1392.1476
$ws.Range("J141").Value = 7338.4546
$ws.Range("K141").Value = 4176.4428
$ws.Range("L141").Value = 22015.3638
$ws.Range("M141").Value = 1003.5572
$ws.Range("N141").Value = -32375.3638

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15811.25
$ws.Range("I32").Value = 2352.8193
$ws.Range("K32").Value = 2352.8193
$ws.Range("M32").Value = -2065.8193

# Row 61
$ws.Range("H61").Value = 2338.6
$ws.Range("I61").Value = 1800.6086
$ws.Range("K61").Value = 1800.6086
$ws.Range("M61").Value = -1588.6086

# Row 97
$ws.Range("H97").Value = 6346.8237
$ws.Range("I97").Value = 6346.8237
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 6346.8237
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -5850.8237
$ws.Range("N97").ClearContents()

# Row 132
$ws.Range("H132").Value = 2048.623
$ws.Range("I132").Value = 1930.6342
$ws.Range("K132").Value = 5791.902599999999
$ws.Range("M132").Value = -3261.902599999999

# Row 136
$ws.Range("H136").Value = 2338.6
$ws.Range("I136").Value = 1800.6086
$ws.Range("K136").Value = 5401.825800000001
$ws.Range("M136").Value = -2851.825800000001

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 21279304
$ws.Range("I134").Value = 28573292
$ws.Range("J134").Value = 5176.3335
$ws.Range("K134").Value = 85719876
$ws.Range("L134").Value = 15529.0005
$ws.Range("M134").Value = -85717341
$ws.Range("N134").Value = -20599.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2112
$ws.Range("I31").Value = 1256.0358
$ws.Range("J31").Value = 4109.25
$ws.Range("K31").Value = 1256.0358
$ws.Range("L31").Value = 4109.25
$ws.Range("M31").Value = -961.0358000000001
$ws.Range("N31").Value = -4699.25

# Row 34
$ws.Range("H34").Value = 2112
$ws.Range("I34").Value = 1256.0358
$ws.Range("J34").Value = 4109.25
$ws.Range("K34").Value = 1256.0358
$ws.Range("L34").Value = 4109.25
$ws.Range("M34").Value = -1054.0358
$ws.Range("N34").Value = -4513.25

# Row 58
$ws.Range("H58").Value = 2901.8333
$ws.Range("I58").Value = 2081.4
$ws.Range("J58").Value = 3927.375
$ws.Range("K58").Value = 2081.4
$ws.Range("L58").Value = 3927.375
$ws.Range("M58").Value = -1878.4
$ws.Range("N58").Value = -4333.375

# Row 132
$ws.Range("H132").Value = 1393.4286
$ws.Range("I132").Value = 1259.2858
$ws.Range("J132").Value = 2332.4285
$ws.Range("K132").Value = 3777.8574
$ws.Range("L132").Value = 6997.2855
$ws.Range("M132").Value = -1247.8574
$ws.Range("N132").Value = -12057.2855

# Row 134
$ws.Range("H134").Value = 1659.1566
$ws.Range("I134").Value = 1122.7333
$ws.Range("J134").Value = 3058.5217
$ws.Range("K134").Value = 3368.199900000001
$ws.Range("L134").Value = 9175.5651
$ws.Range("M134").Value = -833.1999000000005
$ws.Range("N134").Value = -14245.5651

# Row 136
$ws.Range("H136").Value = 2901.8333
$ws.Range("I136").Value = 2081.4
$ws.Range("J136").Value = 3927.375
$ws.Range("K136").Value = 6244.200000000001
$ws.Range("L136").Value = 11782.125
$ws.Range("M136").Value = -3694.200000000001
$ws.Range("N136").Value = -16882.125

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1293
$ws.Range("I5").Value = 939.6667
$ws.Range("K5").Value = 2819.0001
$ws.Range("M5").Value = -2707.0001

# Row 130
$ws.Range("H130").Value = 1000
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# Row 131
$ws.Range("H131").Value = 1601.8108
$ws.Range("J131").Value = 1689.3235
$ws.Range("L131").Value = 5067.970499999999
$ws.Range("N131").Value = -15147.9705

# Row 135
$ws.Range("H135").Value = 1293
$ws.Range("I135").Value = 939.6667
$ws.Range("K135").Value = 8457.0003
$ws.Range("M135").Value = -5922.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 24000
$ws.Range("J34").Value = 24000
$ws.Range("L34").Value = 24000
$ws.Range("N34").Value = -24536

# Row 76
$ws.Range("H76").Value = 24000
$ws.Range("J76").Value = 24000
$ws.Range("L76").Value = 24000
$ws.Range("N76").Value = -24630

# Row 79
$ws.Range("H79").Value = 24000
$ws.Range("J79").Value = 24000
$ws.Range("L79").Value = 24000
$ws.Range("N79").Value = -26184

# Row 113
$ws.Range("H113").Value = 2039.4166
$ws.Range("I113").Value = 1593.9231
$ws.Range("K113").Value = 1593.9231
$ws.Range("M113").Value = 576.0769

# Row 132
$ws.Range("H132").Value = 3251.6667
$ws.Range("I132").Value = 3081.606
$ws.Range("J132").Value = 3719.3333
$ws.Range("K132").Value = 9244.818000000001
$ws.Range("L132").Value = 11157.9999
$ws.Range("M132").Value = -6714.818000000001
$ws.Range("N132").Value = -16217.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4868.154
$ws.Range("I61").Value = 5189.364
$ws.Range("J61").Value = 3101.5
$ws.Range("K61").Value = 5189.364
$ws.Range("L61").Value = 3101.5
$ws.Range("M61").Value = -4987.364
$ws.Range("N61").Value = -3505.5

# Row 100
$ws.Range("H100").Value = 2753.9583
$ws.Range("I100").Value = 2400.5
$ws.Range("J100").Value = 2871.7778
$ws.Range("K100").Value = 2400.5
$ws.Range("L100").Value = 2871.7778
$ws.Range("M100").Value = -1859.5
$ws.Range("N100").Value = -3953.7778

# Row 113
$ws.Range("H113").Value = 4868.154
$ws.Range("I113").Value = 5189.364
$ws.Range("J113").Value = 3101.5
$ws.Range("K113").Value = 5189.364
$ws.Range("L113").Value = 3101.5
$ws.Range("M113").Value = -3019.364
$ws.Range("N113").Value = -7441.5

# Row 132
$ws.Range("H132").Value = 7826.8
$ws.Range("I132").Value = 8489.348
$ws.Range("J132").Value = 6556.9165
$ws.Range("K132").Value = 25468.044
$ws.Range("L132").Value = 19670.7495
$ws.Range("M132").Value = -22938.044
$ws.Range("N132").Value = -24730.7495

$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

# Row 132
$ws.Range("H132").Value = 8773874
$ws.Range("I132").Value = 12196804
$ws.Range("J132").Value = 2617
$ws.Range("K132").Value = 36590412
$ws.Range("L132").Value = 7851
$ws.Range("M132").Value = -36587882
$ws.Range("N132").Value = -12911
